$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.045206739554367
$ws.Cells.Item(2, 4).Value = 1.047739943310706
$ws.Cells.Item(2, 5).Value = 1.052675373424201
$ws.Cells.Item(2, 6).Value = 1.06289267678469
$ws.Cells.Item(2, 9).Value = 1.044393263772516
$ws.Cells.Item(2, 10).Value = 1.0502677111659
$ws.Cells.Item(2, 11).Value = 1.050501737481981
$ws.Cells.Item(2, 12).Value = 1.05542344964904
$ws.Cells.Item(2, 13).Value = 1.065612792522559
$ws.Cells.Item(2, 14).Value = 1.020550180742235

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046293951998988
$ws.Cells.Item(3, 4).Value = 1.04858944930405
$ws.Cells.Item(3, 5).Value = 1.053715903083421
$ws.Cells.Item(3, 6).Value = 1.06418324364314
$ws.Cells.Item(3, 9).Value = 1.044737263987275
$ws.Cells.Item(3, 10).Value = 1.051001811222613
$ws.Cells.Item(3, 11).Value = 1.051163050957941
$ws.Cells.Item(3, 12).Value = 1.056276288301577
$ws.Cells.Item(3, 13).Value = 1.066717069276222
$ws.Cells.Item(3, 14).Value = 1.020800729508543

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.046997246788493
$ws.Cells.Item(4, 4).Value = 1.049138851339466
$ws.Cells.Item(4, 5).Value = 1.054389381715812
$ws.Cells.Item(4, 6).Value = 1.065018985639093
$ws.Cells.Item(4, 9).Value = 1.044958384590329
$ws.Cells.Item(4, 10).Value = 1.051476047858435
$ws.Cells.Item(4, 11).Value = 1.051590037097833
$ws.Cells.Item(4, 12).Value = 1.056827714120082
$ws.Cells.Item(4, 13).Value = 1.067431713460054
$ws.Cells.Item(4, 14).Value = 1.020962434311278

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.047292863860734
$ws.Cells.Item(5, 4).Value = 1.049369751732415
$ws.Cells.Item(5, 5).Value = 1.05467255725176
$ws.Cells.Item(5, 6).Value = 1.065370490438074
$ws.Cells.Item(5, 9).Value = 1.045050991798726
$ws.Cells.Item(5, 10).Value = 1.051675231723966
$ws.Cells.Item(5, 11).Value = 1.051769320146046
$ws.Cells.Item(5, 12).Value = 1.05705943401484
$ws.Cells.Item(5, 13).Value = 1.067732175486878
$ws.Cells.Item(5, 14).Value = 1.021030315206699

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.047342496475098
$ws.Cells.Item(6, 4).Value = 1.049408516924979
$ws.Cells.Item(6, 5).Value = 1.054720106330718
$ws.Cells.Item(6, 6).Value = 1.065429518988784
$ws.Cells.Item(6, 9).Value = 1.045066520325887
$ws.Cells.Item(6, 10).Value = 1.051708664720908
$ws.Cells.Item(6, 11).Value = 1.051799409551207
$ws.Cells.Item(6, 12).Value = 1.057098334983691
$ws.Cells.Item(6, 13).Value = 1.067782625910985
$ws.Cells.Item(6, 14).Value = 1.021041706856003

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.04700119702602
$ws.Cells.Item(7, 4).Value = 1.049141936909839
$ws.Cells.Item(7, 5).Value = 1.054393165343264
$ws.Cells.Item(7, 6).Value = 1.065023681837792
$ws.Cells.Item(7, 9).Value = 1.044959623394885
$ws.Cells.Item(7, 10).Value = 1.051478710090311
$ws.Cells.Item(7, 11).Value = 1.051592433558444
$ws.Cells.Item(7, 12).Value = 1.056830810763063
$ws.Cells.Item(7, 13).Value = 1.06743572814661
$ws.Cells.Item(7, 14).Value = 1.020963341731248

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.045574210926877
$ws.Cells.Item(8, 4).Value = 1.048027096906607
$ws.Cells.Item(8, 5).Value = 1.053026987009093
$ws.Cells.Item(8, 6).Value = 1.063328694127599
$ws.Cells.Item(8, 9).Value = 1.04450982483505
$ws.Cells.Item(8, 10).Value = 1.050515964445542
$ws.Cells.Item(8, 11).Value = 1.050725423725993
$ws.Cells.Item(8, 12).Value = 1.05571175728175
$ws.Cells.Item(8, 13).Value = 1.065985967168878
$ws.Cells.Item(8, 14).Value = 1.020634941230905

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.04305807472065
$ws.Cells.Item(9, 4).Value = 1.046060419947888
$ws.Cells.Item(9, 5).Value = 1.050621010599323
$ws.Cells.Item(9, 6).Value = 1.060346905872601
$ws.Cells.Item(9, 9).Value = 1.043705952413457
$ws.Cells.Item(9, 10).Value = 1.048813529789904
$ws.Cells.Item(9, 11).Value = 1.049190524858221
$ws.Cells.Item(9, 12).Value = 1.053736619009578
$ws.Cells.Item(9, 13).Value = 1.063432054170899
$ws.Cells.Item(9, 14).Value = 1.020053062416614

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041379514375625
$ws.Cells.Item(10, 4).Value = 1.044747825365715
$ws.Cells.Item(10, 5).Value = 1.04901793566044
$ws.Cells.Item(10, 6).Value = 1.058362333906462
$ws.Cells.Item(10, 9).Value = 1.043162448406282
$ws.Cells.Item(10, 10).Value = 1.047674544463745
$ws.Cells.Item(10, 11).Value = 1.048162457411649
$ws.Cells.Item(10, 12).Value = 1.05241765848723
$ws.Cells.Item(10, 13).Value = 1.061729892103281
$ws.Cells.Item(10, 14).Value = 1.019662991978727

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040652396072372
$ws.Cells.Item(11, 4).Value = 1.044179104780985
$ws.Cells.Item(11, 5).Value = 1.048323993561367
$ws.Cells.Item(11, 6).Value = 1.057503754956346
$ws.Cells.Item(11, 9).Value = 1.042925302002701
$ws.Cells.Item(11, 10).Value = 1.047180389130922
$ws.Cells.Item(11, 11).Value = 1.047716149533988
$ws.Cells.Item(11, 12).Value = 1.051846003335379
$ws.Cells.Item(11, 13).Value = 1.060992929272824
$ws.Cells.Item(11, 14).Value = 1.019493575957855

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.040382267018086
$ws.Cells.Item(12, 4).Value = 1.043967802461428
$ws.Cells.Item(12, 5).Value = 1.048066261759222
$ws.Cells.Item(12, 6).Value = 1.057184952195822
$ws.Cells.Item(12, 9).Value = 1.042836943654871
$ws.Cells.Item(12, 10).Value = 1.046996691873676
$ws.Cells.Item(12, 11).Value = 1.047550197976879
$ws.Cells.Item(12, 12).Value = 1.051633583723628
$ws.Cells.Item(12, 13).Value = 1.060719199884361
$ws.Cells.Item(12, 14).Value = 1.019430570116187

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.040440212671512
$ws.Cells.Item(13, 4).Value = 1.044013129982776
$ws.Cells.Item(13, 5).Value = 1.048121544798307
$ws.Cells.Item(13, 6).Value = 1.057253331441032
$ws.Cells.Item(13, 9).Value = 1.042855909117119
$ws.Cells.Item(13, 10).Value = 1.047036102174709
$ws.Cells.Item(13, 11).Value = 1.047585802982856
$ws.Cells.Item(13, 12).Value = 1.051679152149823
$ws.Cells.Item(13, 13).Value = 1.060777915271268
$ws.Cells.Item(13, 14).Value = 1.019444088578029

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040630068044692
$ws.Cells.Item(14, 4).Value = 1.044161639572271
$ws.Cells.Item(14, 5).Value = 1.048302688779606
$ws.Cells.Item(14, 6).Value = 1.057477400351303
$ws.Cells.Item(14, 9).Value = 1.042918003816966
$ws.Cells.Item(14, 10).Value = 1.047165207639773
$ws.Cells.Item(14, 11).Value = 1.047702435457218
$ws.Cells.Item(14, 12).Value = 1.051828446330046
$ws.Cells.Item(14, 13).Value = 1.060970302490661
$ws.Cells.Item(14, 14).Value = 1.019488369450662

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040747038212563
$ws.Cells.Item(15, 4).Value = 1.044253134043947
$ws.Cells.Item(15, 5).Value = 1.048414301407793
$ws.Cells.Item(15, 6).Value = 1.057615471355412
$ws.Cells.Item(15, 9).Value = 1.042956226408992
$ws.Cells.Item(15, 10).Value = 1.04724473441642
$ws.Cells.Item(15, 11).Value = 1.047774273643415
$ws.Cells.Item(15, 12).Value = 1.051920420586109
$ws.Cells.Item(15, 13).Value = 1.061088840080887
$ws.Cells.Item(15, 14).Value = 1.019515642125841

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041427764590702
$ws.Cells.Item(16, 4).Value = 1.044785561898453
$ws.Cells.Item(16, 5).Value = 1.049063994452212
$ws.Cells.Item(16, 6).Value = 1.05841933071231
$ws.Cells.Item(16, 9).Value = 1.043178148960976
$ws.Cells.Item(16, 10).Value = 1.047707319492427
$ws.Cells.Item(16, 11).Value = 1.048192053169799
$ws.Cells.Item(16, 12).Value = 1.052455585982146
$ws.Cells.Item(16, 13).Value = 1.061778803571326
$ws.Cells.Item(16, 14).Value = 1.019674224744315

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.04185468755737
$ws.Cells.Item(17, 4).Value = 1.045119443630931
$ws.Cells.Item(17, 5).Value = 1.049471582664531
$ws.Cells.Item(17, 6).Value = 1.058923771163182
$ws.Cells.Item(17, 9).Value = 1.043316871464991
$ws.Cells.Item(17, 10).Value = 1.047997227486586
$ws.Cells.Item(17, 11).Value = 1.048453807567172
$ws.Cells.Item(17, 12).Value = 1.052791136866147
$ws.Cells.Item(17, 13).Value = 1.062211621176142
$ws.Cells.Item(17, 14).Value = 1.019773562013558

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.042103676548579
$ws.Cells.Item(18, 4).Value = 1.045314156631049
$ws.Cells.Item(18, 5).Value = 1.049709341319932
$ws.Cells.Item(18, 6).Value = 1.059218075579955
$ws.Cells.Item(18, 9).Value = 1.043397611765371
$ws.Cells.Item(18, 10).Value = 1.048166232641035
$ws.Cells.Item(18, 11).Value = 1.048606373756444
$ws.Cells.Item(18, 12).Value = 1.052986806306676
$ws.Cells.Item(18, 13).Value = 1.062464084609166
$ws.Cells.Item(18, 14).Value = 1.019831454278427

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.042188570657927
$ws.Cells.Item(19, 4).Value = 1.045380542878882
$ws.Cells.Item(19, 5).Value = 1.049790414184472
$ws.Cells.Item(19, 6).Value = 1.059318438182337
$ws.Cells.Item(19, 9).Value = 1.043425112586299
$ws.Cells.Item(19, 10).Value = 1.048223843230366
$ws.Cells.Item(19, 11).Value = 1.048658376118548
$ws.Cells.Item(19, 12).Value = 1.053053515783025
$ws.Cells.Item(19, 13).Value = 1.062550169553524
$ws.Cells.Item(19, 14).Value = 1.019851185655148

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.041808885663695
$ws.Cells.Item(20, 4).Value = 1.045083624858285
$ws.Cells.Item(20, 5).Value = 1.049427850311506
$ws.Cells.Item(20, 6).Value = 1.05886964195403
$ws.Cells.Item(20, 9).Value = 1.043302005872786
$ws.Cells.Item(20, 10).Value = 1.047966132754815
$ws.Cells.Item(20, 11).Value = 1.048425735253415
$ws.Cells.Item(20, 12).Value = 1.052755140786642
$ws.Cells.Item(20, 13).Value = 1.062165183099369
$ws.Cells.Item(20, 14).Value = 1.019762909183858

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.04057416162022
$ws.Cells.Item(21, 4).Value = 1.04411790869764
$ws.Cells.Item(21, 5).Value = 1.048249345592985
$ws.Cells.Item(21, 6).Value = 1.057411414571351
$ws.Cells.Item(21, 9).Value = 1.042899725976206
$ws.Cells.Item(21, 10).Value = 1.047127193332802
$ws.Cells.Item(21, 11).Value = 1.047668094878364
$ws.Cells.Item(21, 12).Value = 1.051784485169682
$ws.Cells.Item(21, 13).Value = 1.060913648944758
$ws.Cells.Item(21, 14).Value = 1.019475331969113

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039797580954727
$ws.Cells.Item(22, 4).Value = 1.043510411036074
$ws.Cells.Item(22, 5).Value = 1.047508541373493
$ws.Cells.Item(22, 6).Value = 1.056495214568705
$ws.Cells.Item(22, 9).Value = 1.042645224988807
$ws.Cells.Item(22, 10).Value = 1.046598873989195
$ws.Cells.Item(22, 11).Value = 1.047190735489253
$ws.Cells.Item(22, 12).Value = 1.051173723888449
$ws.Cells.Item(22, 13).Value = 1.060126825184279
$ws.Cells.Item(22, 14).Value = 1.019294074111159

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.040209286194477
$ws.Cells.Item(23, 4).Value = 1.043832486871297
$ws.Cells.Item(23, 5).Value = 1.047901240020433
$ws.Cells.Item(23, 6).Value = 1.056980848867402
$ws.Cells.Item(23, 9).Value = 1.042780289867602
$ws.Cells.Item(23, 10).Value = 1.046879026369832
$ws.Cells.Item(23, 11).Value = 1.047443887760636
$ws.Cells.Item(23, 12).Value = 1.051497544975125
$ws.Cells.Item(23, 13).Value = 1.060543929574195
$ws.Cells.Item(23, 14).Value = 1.01939020470984

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.041829581653904
$ws.Cells.Item(24, 4).Value = 1.045099809926452
$ws.Cells.Item(24, 5).Value = 1.049447611019322
$ws.Cells.Item(24, 6).Value = 1.058894100388177
$ws.Cells.Item(24, 9).Value = 1.043308723532088
$ws.Cells.Item(24, 10).Value = 1.047980183413597
$ws.Cells.Item(24, 11).Value = 1.048438420265545
$ws.Cells.Item(24, 12).Value = 1.052771406025035
$ws.Cells.Item(24, 13).Value = 1.062186166441486
$ws.Cells.Item(24, 14).Value = 1.019767722891881

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043708752287248
$ws.Cells.Item(25, 4).Value = 1.046569112847883
$ws.Cells.Item(25, 5).Value = 1.05124285033374
$ws.Cells.Item(25, 6).Value = 1.061117184957667
$ws.Cells.Item(25, 9).Value = 1.043915109457342
$ws.Cells.Item(25, 10).Value = 1.049254358171599
$ws.Cells.Item(25, 11).Value = 1.049588178082641
$ws.Cells.Item(25, 12).Value = 1.054247625112327
$ws.Cells.Item(25, 13).Value = 1.064092218903546
$ws.Cells.Item(25, 14).Value = 1.020203870985841
